$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Text columns (Coin name, Link) - plain string values, never numeric-looking
$textEdits = @(
    @('B7', 'MXToken'),
    @('C7', 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'),
    @('B8', 'FTXToken'),
    @('C8', 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'),
    @('B9', 'One'),
    @('C9', 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'),
    @('B10', 'WazirX'),
    @('C10', 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'),
    @('B11', 'MandalaExchangeToken'),
    @('C11', 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'),
    @('B12', 'BitrueCoin'),
    @('C12', 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'),
    @('B13', 'BitMartToken'),
    @('C13', 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'),
    @('B14', 'BitForexToken'),
    @('C14', 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'),
    @('B15', 'TigerCash'),
    @('C15', 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'),
    @('B16', 'UpBots'),
    @('C16', 'https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt'),
    @('B17', 'LEO'),
    @('C17', 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'),
    @('B18', 'GateToken'),
    @('C18', 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'),
    @('B26', 'BitKan'),
    @('C26', 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'),
    @('B27', 'HotbitToken'),
    @('C27', 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'),
)

# Numeric-looking text columns (Index, Price, Volume%, Hora) - must stay as Text type
# like the original inline-string cells, so force text entry mode.
$numericLikeEdits = @(
    @('D2', '247.27'),
    @('E2', '1.01%'),
    @('G2', '22'),
    @('D3', '29.41'),
    @('E3', '7.30%'),
    @('G3', '22'),
    @('D4', '5.196'),
    @('E4', '1.56%'),
    @('G4', '22'),
    @('D5', '0.05737'),
    @('E5', '0.90%'),
    @('G5', '22'),
    @('D6', '6.559'),
    @('E6', '0.66%'),
    @('G6', '22'),
    @('D7', '0.8584'),
    @('E7', '4.62%'),
    @('G7', '22'),
    @('D8', '0.8700'),
    @('E8', '2.00%'),
    @('G8', '22'),
    @('D9', '0.01031'),
    @('E9', '1,613.01%'),
    @('G9', '22'),
    @('D10', '0.1367'),
    @('E10', '2.44%'),
    @('G10', '22'),
    @('D11', '0.07066'),
    @('E11', '1.81%'),
    @('G11', '22'),
    @('D12', '0.03071'),
    @('E12', '6.75%'),
    @('G12', '22'),
    @('D13', '0.09385'),
    @('E13', '-0.10%'),
    @('G13', '22'),
    @('D14', '0.001550'),
    @('E14', '1.50%'),
    @('G14', '22'),
    @('D15', '0.006091'),
    @('E15', '-2.00%'),
    @('G15', '22'),
    @('D16', '0.007489'),
    @('E16', '5,224.81%'),
    @('G16', '22'),
    @('D17', '3.494'),
    @('E17', '-0.53%'),
    @('G17', '22'),
    @('D18', '3.100'),
    @('E18', '3.00%'),
    @('G18', '22'),
    @('D19', '2.279'),
    @('E19', '-1.77%'),
    @('G19', '22'),
    @('D20', '0.3184'),
    @('E20', '1.12%'),
    @('G20', '22'),
    @('D21', '0.03322'),
    @('E21', '3.27%'),
    @('G21', '22'),
    @('D22', '0.1290'),
    @('E22', '1.29%'),
    @('G22', '22'),
    @('D23', '3.465'),
    @('E23', '-2.53%'),
    @('G23', '22'),
    @('D24', '0.04136'),
    @('E24', '2.88%'),
    @('G24', '22'),
    @('D25', '0.1380'),
    @('G25', '22'),
    @('D26', '0.001226'),
    @('E26', '0.91%'),
    @('G26', '22'),
    @('D27', '0.004996'),
    @('E27', '11.54%'),
    @('G27', '22'),
    @('D28', '0.0001210'),
    @('E28', '2.51%'),
    @('G28', '22'),
    @('G29', '22'),
    @('G30', '22'),
    @('G31', '22'),
    @('G32', '22'),
    @('G33', '22'),
    @('G34', '22'),
    @('G35', '22'),
    @('G36', '22'),
    @('G37', '22'),
    @('G38', '22'),
    @('G39', '22'),
    @('D40', '0.03758'),
    @('E40', '1.12%'),
    @('G40', '22'),
    @('D41', '0.005749'),
    @('E41', '-3.72%'),
    @('G41', '22'),
    @('D42', '0.1071'),
    @('E42', '1.18%'),
    @('G42', '22'),
    @('D43', '0.002428'),
    @('E43', '2.20%'),
    @('G43', '22'),
    @('D44', '0.009472'),
    @('E44', '-2.48%'),
    @('G44', '22'),
    @('D45', '0.00005260'),
    @('E45', '2.97%'),
    @('G45', '22'),
    @('G46', '22'),
    @('D47', '0.05699'),
    @('E47', '-43.58%'),
    @('G47', '22'),
    @('D48', '0.002277'),
    @('E48', '-9.35%'),
    @('G48', '22'),
    @('D49', '0.00002100'),
    @('G49', '22'),
    @('D50', '0.0002000'),
    @('G50', '22'),
    @('G51', '22'),
)

foreach ($edit in $textEdits) {
    $ws.Range($edit[0]).Value = $edit[1]
}

foreach ($edit in $numericLikeEdits) {
    $cell = $ws.Range($edit[0])
    $cell.NumberFormat = "@"
    $cell.Value = $edit[1]
    $cell.Style = "Normal"
}